# Add a new "month calendar icon" test treatment (rows 32-34) to the
# "Web Parameters" sheet, mirroring the existing "year" icon treatment
# (treatment_id 16, rows 29-31) but with month-span dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Web Parameters")

# --- New row 32 ---------------------------------------------------------
$ws.Range("A32").Value = 17
$ws.Range("B32").Value = 1
$ws.Range("C32").Value = "calendarIcon"
$ws.Range("D32").Value = "none"
$ws.Range("E32").Value = "none"
$ws.Range("F32").Value = 300
$ws.Range("H32").Value = 44593
$ws.Range("I32").Value = 700
$ws.Range("K32").Value = 44614
$ws.Range("L32").Value = 1100
$ws.Range("N32").Value = 100
$ws.Range("O32").Value = 100
$ws.Range("T32").Value = 8
$ws.Range("U32").Value = 8
$ws.Range("V32").Value = "Calendar month view with icon and no interaction."

# --- New row 33 ---------------------------------------------------------
$ws.Range("A33").Value = 17
$ws.Range("B33").Value = 2
$ws.Range("C33").Value = "calendarIcon"
$ws.Range("D33").Value = "none"
$ws.Range("E33").Value = "none"
$ws.Range("F33").Value = 500
$ws.Range("H33").Value = 44621
$ws.Range("I33").Value = 800
$ws.Range("K33").Value = 44632
$ws.Range("L33").Value = 1100
$ws.Range("N33").Value = 100
$ws.Range("O33").Value = 100
$ws.Range("T33").Value = 8
$ws.Range("U33").Value = 8
$ws.Range("V33").Value = "Calendar month view with icon and no interaction."

# --- New row 34 ---------------------------------------------------------
$ws.Range("A34").Value = 17
$ws.Range("B34").Value = 3
$ws.Range("C34").Value = "calendarIcon"
$ws.Range("D34").Value = "none"
$ws.Range("E34").Value = "none"
$ws.Range("F34").Value = 300
$ws.Range("H34").Value = 44652
$ws.Range("I34").Value = 1000
$ws.Range("K34").Value = 44666
$ws.Range("L34").Value = 1100
$ws.Range("N34").Value = 100
$ws.Range("O34").Value = 100
$ws.Range("T34").Value = 8
$ws.Range("U34").Value = 8
$ws.Range("V34").Value = "Calendar month view with icon and no interaction."

# --- Restore the view: scroll back to column A and select A31 ----------
$ws.Activate() | Out-Null
$ws.Range("A31").Select() | Out-Null
